$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New values (file #1 renamed ccf55d22-... -> 25f65b42-...,
#             file #2 renamed e35b6dea-...  -> ffff0088d0ca-...)
# ---------------------------------------------------------------------------
$newFile1 = "25f65b42-0d0c-4056-b77d-c73e5608bc99.md"
$newPath1 = "e2e\25f65b42-0d0c-4056-b77d-c73e5608bc99.md"
$newFile2 = "ffff0088d0ca-cb3b-41e2-a7d9-fb77da1f136b.md"
$newPath2 = "e2e\ffff0088d0ca-cb3b-41e2-a7d9-fb77da1f136b.md"

$oldFile1 = "ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.md"
$oldFile2 = "e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.md"

$dateOverview = "2016-08-22 15:13:18"

$zhcnXlf1    = "25f65b42-0d0c-4056-b77d-c73e5608bc99.09afa1eb048dc65acb25fa94f28087177b79b34b.zh-cn.xlf"
$zhcnDateIn  = "2016-08-22 15:13:11"
$zhcnDateOut = "2016-08-22 15:13:36"

$dedeXlf1    = "25f65b42-0d0c-4056-b77d-c73e5608bc99.09afa1eb048dc65acb25fa94f28087177b79b34b.de-de.xlf"
$dedeDateOut = "2016-08-22 15:13:43"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("G2").Value = $dateOverview
$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("G3").Value = $dateOverview

# Hyperlinks.Delete() on a range removes every hyperlink on the sheet in this
# engine, so recreate all of them (targets are unchanged, only the display
# text changes).
$wsOverview.Range("B2").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile1",
    "",
    "",
    $newPath1
)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile2",
    "",
    "",
    $newPath2
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $zhcnXlf1
$wsZhCn.Range("H2").Value = $zhcnDateIn
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $zhcnXlf1
$wsZhCn.Range("K2").Value = $zhcnDateOut

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $zhcnXlf1
$wsZhCn.Range("H3").Value = $zhcnDateIn
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $zhcnXlf1
$wsZhCn.Range("K3").Value = $zhcnDateOut

$wsZhCn.Range("A2").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile1",
    "",
    "",
    $newFile1
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c6836b8fd7f80834c7261cebec0eccf4b705e15b/e2e/$oldFile1",
    "",
    "",
    $newFile1
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile2",
    "",
    "",
    $newFile2
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c6836b8fd7f80834c7261cebec0eccf4b705e15b/e2e/$oldFile2",
    "",
    "",
    $newFile2
)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $dedeXlf1
$wsDeDe.Range("H2").Value = $dateOverview
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $dedeXlf1
$wsDeDe.Range("K2").Value = $dedeDateOut

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $dedeXlf1
$wsDeDe.Range("H3").Value = $dateOverview
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $dedeXlf1
$wsDeDe.Range("K3").Value = $dedeDateOut

$wsDeDe.Range("A2").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile1",
    "",
    "",
    $newFile1
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/39e33e9c4ff97edc5c66fcb4466fbed3544cbb35/e2e/$oldFile1",
    "",
    "",
    $newFile1
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/$oldFile2",
    "",
    "",
    $newFile2
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/39e33e9c4ff97edc5c66fcb4466fbed3544cbb35/e2e/$oldFile2",
    "",
    "",
    $newFile2
)
